$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the existing header formatting
# (bold font, thin border, centered/top-aligned) by copying the format
# from the last existing header cell (G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Corresponding data value for row 2
$ws.Range("H2").Value = 0
